$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''29.102.18'
$ws.Cells.Item(2, 5).Value = '  +0.06%  '
$ws.Cells.Item(3, 4).Value = '''1.821.24'
$ws.Cells.Item(3, 5).Value = '  -0.63%  '
$ws.Cells.Item(4, 4).Value = '''0.9979'
$ws.Cells.Item(5, 4).Value = '''241.43'
$ws.Cells.Item(5, 5).Value = '  -0.80%  '
$ws.Cells.Item(6, 4).Value = '''0.6157'
$ws.Cells.Item(6, 5).Value = '  -2.02%  '
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$ws.Cells.Item(8, 4).Value = '''0.07329'
$ws.Cells.Item(8, 5).Value = '  -2.27%  '
$ws.Cells.Item(9, 4).Value = '''0.2885'
$ws.Cells.Item(9, 5).Value = '  -1.23%  '
$ws.Cells.Item(10, 4).Value = '''22.83'
$ws.Cells.Item(10, 5).Value = '  -1.85%  '
$ws.Cells.Item(11, 4).Value = '''0.07662'
$ws.Cells.Item(11, 5).Value = '  -0.30%  '
$ws.Cells.Item(12, 4).Value = '''1.820.56'
$ws.Cells.Item(13, 4).Value = '''4.943'
$ws.Cells.Item(13, 5).Value = '  -1.28%  '
$ws.Cells.Item(14, 4).Value = '''0.6594'
$ws.Cells.Item(14, 5).Value = '  -1.14%  '
$ws.Cells.Item(15, 4).Value = '''81.52'
$ws.Cells.Item(15, 5).Value = '  -1.48%  '
$ws.Cells.Item(16, 4).Value = '''0.000008964'
$ws.Cells.Item(16, 5).Value = '  -4.41%  '
$ws.Cells.Item(17, 4).Value = '''5.834'
$ws.Cells.Item(17, 5).Value = '  -2.45%  '
$ws.Cells.Item(18, 4).Value = '''29.070.17'
$ws.Cells.Item(18, 5).Value = '  -0.08%  '
$ws.Cells.Item(19, 4).Value = '''2.061.95'
$ws.Cells.Item(19, 5).Value = '  -0.89%  '
$ws.Cells.Item(20, 4).Value = '''237.60'
$ws.Cells.Item(20, 5).Value = '  +6.43%  '
$ws.Cells.Item(21, 4).Value = '''12.43'
$ws.Cells.Item(21, 5).Value = '  -1.22%  '
$ws.Cells.Item(22, 5).Value = '  -0.27%  '
$ws.Cells.Item(23, 4).Value = '''7.103'
$ws.Cells.Item(23, 5).Value = '  +0.04%  '
$ws.Cells.Item(24, 4).Value = '''1.001'
$ws.Cells.Item(24, 5).Value = '  +0.04%  '
$ws.Cells.Item(25, 5).Value = '  -1.35%  '
$ws.Cells.Item(26, 4).Value = '''0.1404'
$ws.Cells.Item(27, 4).Value = '''8.425'
$ws.Cells.Item(27, 5).Value = '  -0.78%  '
$ws.Cells.Item(28, 4).Value = '''17.57'
$ws.Cells.Item(28, 5).Value = '  -1.75%  '
$ws.Cells.Item(29, 4).Value = '''1.476'
$ws.Cells.Item(29, 5).Value = '  -1.45%  '
$ws.Cells.Item(30, 4).Value = '''0.05570'
$ws.Cells.Item(30, 5).Value = '  -1.76%  '
$ws.Cells.Item(31, 4).Value = '''4.088'
$ws.Cells.Item(31, 5).Value = '  +0.09%  '
$ws.Cells.Item(32, 4).Value = '''4.093'
$ws.Cells.Item(32, 5).Value = '  -1.42%  '
$ws.Cells.Item(33, 4).Value = '''1.206'
$ws.Cells.Item(33, 5).Value = '  +0.12%  '
$ws.Cells.Item(34, 4).Value = '''0.7340'
$ws.Cells.Item(34, 5).Value = '  -1.05%  '
$ws.Cells.Item(35, 5).Value = '  -1.78%  '
$ws.Cells.Item(36, 4).Value = '''1.128'
$ws.Cells.Item(36, 5).Value = '  -1.05%  '
$ws.Cells.Item(37, 4).Value = '''2.616'
$ws.Cells.Item(37, 5).Value = '  -2.01%  '
$ws.Cells.Item(38, 4).Value = '''2.828'
$ws.Cells.Item(38, 5).Value = '  +2.40%  '
$ws.Cells.Item(39, 4).Value = '''1.208.49'
$ws.Cells.Item(39, 5).Value = '  -1.09%  '
$ws.Cells.Item(40, 4).Value = '''0.01755'
$ws.Cells.Item(40, 5).Value = '  -1.29%  '
$ws.Cells.Item(41, 4).Value = '''6.375'
$ws.Cells.Item(41, 5).Value = '  -2.40%  '
$ws.Cells.Item(42, 4).Value = '''0.8918'
$ws.Cells.Item(42, 5).Value = '  -0.15%  '
$ws.Cells.Item(43, 5).Value = '  +0.02%  '
$ws.Cells.Item(44, 5).Value = '  -1.19%  '
$ws.Cells.Item(45, 4).Value = '''1.969.02'
$ws.Cells.Item(45, 5).Value = '  -0.61%  '
$ws.Cells.Item(46, 4).Value = '''64.43'
$ws.Cells.Item(46, 5).Value = '  -2.20%  '
$ws.Cells.Item(47, 4).Value = '''0.00000000120'
$ws.Cells.Item(47, 5).Value = '  -2.69%  '
$ws.Cells.Item(48, 4).Value = '''0.5085'
$ws.Cells.Item(48, 5).Value = '  -0.05%  '
$ws.Cells.Item(49, 2).Value = 'TheSandbox'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(49, 4).Value = '''0.3987'
$ws.Cells.Item(49, 5).Value = '  -2.17%  '
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).Value = '''9.012'
$ws.Cells.Item(50, 5).Value = '  +0.02%  '
$ws.Cells.Item(51, 4).Value = '''0.05754'
$ws.Cells.Item(51, 5).Value = '  -1.11%  '
